$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New day of data: 5/12 (2016-05-12, serial 42502) run report
$ws.Range("A29").Value = 42502
$ws.Range("B29").Value = 141
$ws.Range("C29").Value = 134
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 134
$ws.Range("G29").Value = 0.95
$ws.Range("G29").NumberFormat = "0.0%"
$ws.Range("G29").HorizontalAlignment = -4108
$ws.Range("H29").Value = 44.467661691188411
$ws.Range("I29").Value = 34.116666658082977
$ws.Range("J29").Value = 114.299999991199

$ws.Range("L21").Select() | Out-Null
